# "Forgot to add phone number to Aarhus Havn (otherwise no contact channels are displayed)"
#
# Adds the missing "telefon" (phone) values for the two organisation units
# (Aarhus Havn / HAVNEN and Radhuset / ARHUS) on the "organisationenhed"
# sheet, in column L (header L1 = "telefon").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("organisationenhed")
$ws.Activate()

# Row 2 -> Aarhus Havn
$ws.Range("L2").Value = 12345678
# Row 3 -> Radhuset
$ws.Range("L3").Value = 87654321

# Restore/refresh the on-screen selection the way it ends up after entering
# the two values: cursor moved onto the header area first, then left on the
# last-edited cell in the frozen data pane.
$ws.Range("I1").Select()
$ws.Range("L3").Select()
